$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-10 Sunday" "2024-03-11 Monday"

Replace-Text "162×8=" "892×4="
Replace-Text "967×8=" "589×6="
Replace-Text "521×8=" "441×8="
Replace-Text "652×8=" "206×2="
Replace-Text "116×2=" "344×4="
Replace-Text "724×7=" "233×4="
Replace-Text "577×6=" "736×9="
Replace-Text "726×6=" "481×5="
Replace-Text "113×2=" "392×6="
Replace-Text "538×7=" "268×8="
Replace-Text "885×8=" "318×4="
Replace-Text "554×8=" "479×4="
Replace-Text "207×2=" "434×5="
Replace-Text "433×3=" "793×5="
Replace-Text "140×2=" "499×8="
Replace-Text "527×5=" "892×4="
Replace-Text "641×9=" "586×3="
Replace-Text "738×5=" "158×3="
Replace-Text "806×4=" "228×4="
Replace-Text "264×7=" "680×8="
Replace-Text "267×3=" "317×8="
Replace-Text "601×5=" "264×6="
Replace-Text "887×3=" "452×5="
Replace-Text "287×5=" "504×2="
Replace-Text "439×2=" "994×7="
